$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update reference label in D2 (cell holding the citation key used for the Graham 2018 estimate)
$ws.Range("D2").Value = "graham2018budget"
